$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from its old spot (right before the page
#    break that follows "Goodbye!") to right after "---=== IPC Temperature "
#    (i.e. right before the word that is about to become "Calculator").
#    Bookmarks.Add on an existing bookmark name re-seats it (removes the old
#    one and creates the new one), so we don't need a separate Delete call.
# ---------------------------------------------------------------------------
$findAnalyzer = $d.Content
$findAnalyzer.Find.Execute("Analyzer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$analyzerStart = $findAnalyzer.Start
$analyzerEnd = $findAnalyzer.End
$d.Bookmarks.Add("_GoBack", $d.Range($analyzerStart, $analyzerStart))

# ---------------------------------------------------------------------------
# 2) "Analyzer" -> "Calculator"
#    A plain Find/Replace on this run causes the interpreter's run-merge
#    normalisation pass to fuse it together with the neighbouring runs that
#    happen to share identical run formatting ("---=== IPC Temperature ",
#    " V2.0", " ", "===---"). Those runs must stay distinct (per the target
#    XML), so after the textual replace we re-stamp the internal boundaries
#    by round-tripping each sub-range's FormattedText, which forces the
#    engine to re-split the run without touching its text/formatting.
# ---------------------------------------------------------------------------
$findAnalyzer2 = $d.Content
$findAnalyzer2.Find.Execute("Analyzer", $true, $false, $false, $false, $false, $true, 1, $false, "Calculator", 2)

$findTitle = $d.Content
$findTitle.Find.Execute("---=== IPC Temperature Calculator V2.0 ===---", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$titleBase = $findTitle.Start

# Offsets (relative to $titleBase) of the five original runs:
#   "---=== IPC Temperature " | "Calculator" | " V2.0" | " " | "===---"
$titleBounds = @(0, 23, 33, 38, 39, 45)
for ($i = 1; $i -lt $titleBounds.Length - 2; $i++) {
    $subStart = $titleBase + $titleBounds[$i]
    $subEnd = $titleBase + $titleBounds[$i + 1]
    $subRng = $d.Range($subStart, $subEnd)
    $subFmt = $subRng.FormattedText
    $subRng.FormattedText = $subFmt
}

# ---------------------------------------------------------------------------
# 3) "Please enter the number of days between 3 and 10, inclusive:" becomes
#    three runs: "Please enter the number of days" + "," + " between 3 and
#    10, inclusive:" (all sharing the same Courier/red/yellow-highlight
#    formatting). We copy an existing comma that already carries the exact
#    target run formatting (the one found a little further along, in
#    "...10, inclusive:", which is part of the very same original run) and
#    paste it at the insertion point - this creates a genuinely separate
#    run (rather than a Font-property edit, whose results don't always
#    round-trip eastAsia/cs through this runtime) while preserving every
#    rFonts attribute exactly.
# ---------------------------------------------------------------------------
$findDays = $d.Content
$findDays.Find.Execute("Please enter the number of days between 3 and 10, inclusive:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$daysStart = $findDays.Start

$insertPos = $daysStart + [int]"Please enter the number of days".Length

# Locate the comma inside "...3 and 10, inclusive:" *relative to this same
# match* (rather than re-running Find from the top of the document, which
# would find the unrelated, differently-formatted "...3 and 10, inclusive."
# earlier in the body text).
$relCommaOffset = [int]"Please enter the number of days between 3 and 10".Length
$goodCommaStart = $daysStart + $relCommaOffset
$goodComma = $d.Range($goodCommaStart, $goodCommaStart + 1)
$goodComma.Copy()

$pasteTarget = $d.Range($insertPos, $insertPos)
$pasteTarget.Paste()

Write-Output "done"
